# Apply industry-specific content corrections to the IT Requirements
# Traceability Matrix "Open Action Items" sheet.
#
# Rows 8-17 reuse the same four placeholder strings; this patches them with
# the correct IT / Cloud Infrastructure Migration terminology while leaving
# unrelated owners (Ethics Committee, Compliance Officers) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Open Action Items")

# Column E (Owner) replacements, keyed by row number.
$ownerMap = @{
    8  = "Chief Technology Officer"
    9  = "IT Managers"
    10 = "DevOps Engineers"
    11 = "System Administrators"
    14 = "Chief Technology Officer"
    15 = "IT Managers"
    16 = "DevOps Engineers"
    17 = "System Administrators"
}

foreach ($row in $ownerMap.Keys) {
    $ws.Range("E$row").Value = $ownerMap[$row]
}

# Columns I (Dependencies) and J (Notes) get the same replacement across
# every data row from 8 to 17.
for ($row = 8; $row -le 17; $row++) {
    $ws.Range("I$row").Value = "Dependent on Cloud Infrastructure Migration milestone completion"
    $ws.Range("J$row").Value = "Critical action for Information Technology success"
}
